$p = $ppt.ActivePresentation

$oldUS = "11/29/18"
$newUS = "11/30/18"
$oldBR = "29/11/2018"
$newBR = "30/11/2018"

function Update-DateShape {
    param($shape)

    if (-not $shape.HasTextFrame) { return }
    $tr = $shape.TextFrame.TextRange
    $t = $tr.Text
    if ($t -eq $oldUS) {
        $tr.Text = $newUS
    } elseif ($t -eq $oldBR) {
        $tr.Text = $newBR
    }
}

# Slide master
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# Every slide layout (custom layout) hanging off the slide master
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# Notes master
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    Update-DateShape $notesMaster.Shapes.Item($i)
}
